$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(11, 8).Value = 278.7857  # H11: 148.83333 -> 278.7857
$ws.Cells.Item(11, 9).Value = 278.7857  # I11: 148.83333 -> 278.7857
$ws.Cells.Item(11, 11).Value = 278.7857  # K11: 148.83333 -> 278.7857
$ws.Cells.Item(11, 13).Value = -138.7857  # M11: -8.833329999999989 -> -138.7857
$ws.Cells.Item(64, 8).Value = 7640.55  # H64: 7645.6113 -> 7640.55
$ws.Cells.Item(64, 9).Value = 4321.1  # I64: 4152.1 -> 4321.1
$ws.Cells.Item(64, 10).Value = 10960  # J64: 12012.5 -> 10960
$ws.Cells.Item(64, 11).Value = 4321.1  # K64: 4152.1 -> 4321.1
$ws.Cells.Item(64, 12).Value = 10960  # L64: 12012.5 -> 10960
$ws.Cells.Item(64, 13).Value = -4073.1  # M64: -3904.1 -> -4073.1
$ws.Cells.Item(64, 14).Value = -11456  # N64: -12508.5 -> -11456
$ws.Cells.Item(67, 8).Value = 7640.55  # H67: 7645.6113 -> 7640.55
$ws.Cells.Item(67, 9).Value = 4321.1  # I67: 4152.1 -> 4321.1
$ws.Cells.Item(67, 10).Value = 10960  # J67: 12012.5 -> 10960
$ws.Cells.Item(67, 11).Value = 4321.1  # K67: 4152.1 -> 4321.1
$ws.Cells.Item(67, 12).Value = 10960  # L67: 12012.5 -> 10960
$ws.Cells.Item(67, 13).Value = -3463.1  # M67: -3294.1 -> -3463.1
$ws.Cells.Item(67, 14).Value = -12676  # N67: -13728.5 -> -12676
$ws.Cells.Item(116, 8).Value = 4398.3  # H116: 4542.6665 -> 4398.3
$ws.Cells.Item(116, 9).Value = 3324.75  # I116: 3400 -> 3324.75
$ws.Cells.Item(116, 11).Value = 3324.75  # K116: 3400 -> 3324.75
$ws.Cells.Item(116, 13).Value = 117.25  # M116: 42 -> 117.25
$ws.Cells.Item(131, 8).Value = 56310.45  # H131: 58232.05 -> 56310.45
$ws.Cells.Item(131, 9).Value = 60345  # I131: 60356.055 -> 60345
$ws.Cells.Item(131, 10).Value = 19999.5  # J131: 20000 -> 19999.5
$ws.Cells.Item(131, 11).Value = 181035  # K131: 181068.165 -> 181035
$ws.Cells.Item(131, 12).Value = 59998.5  # L131: 60000 -> 59998.5
$ws.Cells.Item(131, 13).Value = -175995  # M131: -176028.165 -> -175995
$ws.Cells.Item(131, 14).Value = -70078.5  # N131: -70080 -> -70078.5
$ws.Cells.Item(132, 8).Value = 1227.7646  # H132: 997.875 -> 1227.7646
$ws.Cells.Item(132, 9).Value = 1019.44446  # I132: 765.75 -> 1019.44446
$ws.Cells.Item(132, 11).Value = 3058.33338  # K132: 2297.25 -> 3058.33338
$ws.Cells.Item(132, 13).Value = -528.33338  # M132: 232.75 -> -528.33338
$ws.Cells.Item(137, 8).Value = 927.17645  # H137: 1045.4375 -> 927.17645
$ws.Cells.Item(137, 9).Value = 810.86664  # I137: 937.7143 -> 810.86664
$ws.Cells.Item(137, 11).Value = 2432.59992  # K137: 2813.1429 -> 2432.59992
$ws.Cells.Item(137, 13).Value = 117.4000800000003  # M137: -263.1428999999998 -> 117.4000800000003
$ws.Cells.Item(138, 8).Value = 2575.7932  # H138: 2600.1853 -> 2575.7932
$ws.Cells.Item(138, 9).Value = 2145  # I138: 2205.5264 -> 2145
$ws.Cells.Item(138, 10).Value = 3533.111  # J138: 3537.5 -> 3533.111
$ws.Cells.Item(138, 11).Value = 6435  # K138: 6616.5792 -> 6435
$ws.Cells.Item(138, 12).Value = 10599.333  # L138: 10612.5 -> 10599.333
$ws.Cells.Item(138, 13).Value = -1295  # M138: -1476.5792 -> -1295
$ws.Cells.Item(138, 14).Value = -20879.333  # N138: -20892.5 -> -20879.333
$ws.Cells.Item(141, 8).Value = 3265.3914  # H141: 3258.5833 -> 3265.3914
$ws.Cells.Item(141, 9).Value = 3141.0908  # I141: 3139.3914 -> 3141.0908
$ws.Cells.Item(141, 11).Value = 9423.2724  # K141: 9418.174199999999 -> 9423.2724
$ws.Cells.Item(141, 13).Value = -4243.2724  # M141: -4238.174199999999 -> -4243.2724

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(61, 8).Value = 5215300  # H61: 5215362 -> 5215300
$ws.Cells.Item(61, 9).Value = 7581453.5  # I61: 7581544 -> 7581453.5
$ws.Cells.Item(61, 11).Value = 7581453.5  # K61: 7581544 -> 7581453.5
$ws.Cells.Item(61, 13).Value = -7581241.5  # M61: -7581332 -> -7581241.5
$ws.Cells.Item(102, 8).Value = 3777.389  # H102: 3777.4443 -> 3777.389
$ws.Cells.Item(102, 9).Value = 2799.5334  # I102: 2799.6 -> 2799.5334
$ws.Cells.Item(102, 11).Value = 2799.5334  # K102: 2799.6 -> 2799.5334
$ws.Cells.Item(102, 13).Value = -1177.5334  # M102: -1177.6 -> -1177.5334
$ws.Cells.Item(103, 8).Value = 30000  # H103: 0 -> 30000
$ws.Cells.Item(103, 10).Value = 30000  # J103: 0 -> 30000
$ws.Cells.Item(103, 12).Value = 30000  # L103: 0 -> 30000
$ws.Cells.Item(103, 14).Value = -32344  # N103: None -> -32344
$ws.Cells.Item(122, 8).Value = 2615.0435  # H122: 2360.963 -> 2615.0435
$ws.Cells.Item(122, 9).Value = 1902.9474  # I122: 1728.5217 -> 1902.9474
$ws.Cells.Item(122, 11).Value = 5708.8422  # K122: 5185.5651 -> 5708.8422
$ws.Cells.Item(122, 13).Value = -3258.8422  # M122: -2735.5651 -> -3258.8422
$ws.Cells.Item(136, 8).Value = 5215300  # H136: 5215362 -> 5215300
$ws.Cells.Item(136, 9).Value = 7581453.5  # I136: 7581544 -> 7581453.5
$ws.Cells.Item(136, 11).Value = 22744360.5  # K136: 22744632 -> 22744360.5
$ws.Cells.Item(136, 13).Value = -22741810.5  # M136: -22742082 -> -22741810.5

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(94, 8).Value = 2013  # H94: 1948.25 -> 2013
$ws.Cells.Item(94, 9).Value = 1569.3478  # I94: 1512.2916 -> 1569.3478
$ws.Cells.Item(94, 11).Value = 1569.3478  # K94: 1512.2916 -> 1569.3478
$ws.Cells.Item(94, 13).Value = -1118.3478  # M94: -1061.2916 -> -1118.3478
$ws.Cells.Item(107, 8).Value = 4921.4287  # H107: 3276.4 -> 4921.4287
$ws.Cells.Item(107, 9).Value = 3741.6667  # I107: 2653.2856 -> 3741.6667
$ws.Cells.Item(107, 11).Value = 3741.6667  # K107: 2653.2856 -> 3741.6667
$ws.Cells.Item(107, 13).Value = -1821.6667  # M107: -733.2856000000002 -> -1821.6667
$ws.Cells.Item(119, 8).Value = 0  # H119: 35000 -> 0
$ws.Cells.Item(119, 10).Value = 0  # J119: 35000 -> 0
$ws.Cells.Item(119, 12).Value = 0  # L119: 35000 -> 0
$ws.Cells.Item(119, 14).ClearContents()  # N119: remove (was -44676)
$ws.Cells.Item(122, 8).Value = 0  # H122: 20000 -> 0
$ws.Cells.Item(122, 10).Value = 0  # J122: 20000 -> 0
$ws.Cells.Item(122, 12).Value = 0  # L122: 20000 -> 0
$ws.Cells.Item(122, 14).ClearContents()  # N122: remove (was -29800)
$ws.Cells.Item(134, 8).Value = 3280.625  # H134: 3456.9355 -> 3280.625
$ws.Cells.Item(134, 9).Value = 3278.6553  # I134: 3491.5557 -> 3278.6553
$ws.Cells.Item(134, 10).Value = 3299.6667  # J134: 3223.25 -> 3299.6667
$ws.Cells.Item(134, 11).Value = 9835.965899999999  # K134: 10474.6671 -> 9835.965899999999
$ws.Cells.Item(134, 12).Value = 9899.000100000001  # L134: 9669.75 -> 9899.000100000001
$ws.Cells.Item(134, 13).Value = -7300.965899999999  # M134: -7939.667099999999 -> -7300.965899999999
$ws.Cells.Item(134, 14).Value = -14969.0001  # N134: -14739.75 -> -14969.0001

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(19, 8).Value = 736  # H19: 589.9091 -> 736
$ws.Cells.Item(19, 9).Value = 736  # I19: 589.9091 -> 736
$ws.Cells.Item(19, 11).Value = 736  # K19: 589.9091 -> 736
$ws.Cells.Item(19, 13).Value = -566  # M19: -419.9091 -> -566
$ws.Cells.Item(24, 8).Value = 736  # H24: 589.9091 -> 736
$ws.Cells.Item(24, 9).Value = 736  # I24: 589.9091 -> 736
$ws.Cells.Item(24, 11).Value = 736  # K24: 589.9091 -> 736
$ws.Cells.Item(24, 13).Value = -566  # M24: -419.9091 -> -566
$ws.Cells.Item(31, 8).Value = 3695.372  # H31: 3939.975 -> 3695.372
$ws.Cells.Item(31, 9).Value = 2576.1428  # I31: 2664.3704 -> 2576.1428
$ws.Cells.Item(31, 10).Value = 5784.6  # J31: 6589.3076 -> 5784.6
$ws.Cells.Item(31, 11).Value = 2576.1428  # K31: 2664.3704 -> 2576.1428
$ws.Cells.Item(31, 12).Value = 5784.6  # L31: 6589.3076 -> 5784.6
$ws.Cells.Item(31, 13).Value = -2281.1428  # M31: -2369.3704 -> -2281.1428
$ws.Cells.Item(31, 14).Value = -6374.6  # N31: -7179.3076 -> -6374.6
$ws.Cells.Item(34, 8).Value = 3695.372  # H34: 3939.975 -> 3695.372
$ws.Cells.Item(34, 9).Value = 2576.1428  # I34: 2664.3704 -> 2576.1428
$ws.Cells.Item(34, 10).Value = 5784.6  # J34: 6589.3076 -> 5784.6
$ws.Cells.Item(34, 11).Value = 2576.1428  # K34: 2664.3704 -> 2576.1428
$ws.Cells.Item(34, 12).Value = 5784.6  # L34: 6589.3076 -> 5784.6
$ws.Cells.Item(34, 13).Value = -2374.1428  # M34: -2462.3704 -> -2374.1428
$ws.Cells.Item(34, 14).Value = -6188.6  # N34: -6993.3076 -> -6188.6
$ws.Cells.Item(68, 8).Value = 44320  # H68: 44315 -> 44320
$ws.Cells.Item(68, 10).Value = 44320  # J68: 44315 -> 44320
$ws.Cells.Item(68, 12).Value = 44320  # L68: 44315 -> 44320
$ws.Cells.Item(68, 14).Value = -45818  # N68: -45813 -> -45818
$ws.Cells.Item(71, 8).Value = 44320  # H71: 44315 -> 44320
$ws.Cells.Item(71, 10).Value = 44320  # J71: 44315 -> 44320
$ws.Cells.Item(71, 12).Value = 132960  # L71: 132945 -> 132960
$ws.Cells.Item(71, 14).Value = -140448  # N71: -140433 -> -140448
$ws.Cells.Item(74, 8).Value = 30819.8  # H74: 32023.111 -> 30819.8
$ws.Cells.Item(74, 10).Value = 30819.8  # J74: 32023.111 -> 30819.8
$ws.Cells.Item(74, 12).Value = 30819.8  # L74: 32023.111 -> 30819.8
$ws.Cells.Item(74, 14).Value = -32567.8  # N74: -33771.111 -> -32567.8
$ws.Cells.Item(77, 8).Value = 30819.8  # H77: 32023.111 -> 30819.8
$ws.Cells.Item(77, 10).Value = 30819.8  # J77: 32023.111 -> 30819.8
$ws.Cells.Item(77, 12).Value = 92459.39999999999  # L77: 96069.333 -> 92459.39999999999
$ws.Cells.Item(77, 14).Value = -101195.4  # N77: -104805.333 -> -101195.4
$ws.Cells.Item(92, 8).Value = 63514.145  # H92: 58920.2 -> 63514.145
$ws.Cells.Item(92, 10).Value = 63514.145  # J92: 58920.2 -> 63514.145
$ws.Cells.Item(92, 12).Value = 63514.145  # L92: 58920.2 -> 63514.145
$ws.Cells.Item(92, 14).Value = -68506.14499999999  # N92: -63912.2 -> -68506.14499999999
$ws.Cells.Item(99, 8).Value = 6194.857  # H99: 6579.077 -> 6194.857
$ws.Cells.Item(99, 9).Value = 3986.125  # I99: 4384.143 -> 3986.125
$ws.Cells.Item(99, 11).Value = 3986.125  # K99: 4384.143 -> 3986.125
$ws.Cells.Item(99, 13).Value = -2488.125  # M99: -2886.143 -> -2488.125
$ws.Cells.Item(126, 8).Value = 6194.857  # H126: 6579.077 -> 6194.857
$ws.Cells.Item(126, 9).Value = 3986.125  # I126: 4384.143 -> 3986.125
$ws.Cells.Item(126, 11).Value = 11958.375  # K126: 13152.429 -> 11958.375
$ws.Cells.Item(126, 13).Value = -9488.375  # M126: -10682.429 -> -9488.375
$ws.Cells.Item(132, 8).Value = 4743.875  # H132: 2562.7144 -> 4743.875
$ws.Cells.Item(132, 9).Value = 7708  # I132: 1556 -> 7708
$ws.Cells.Item(132, 11).Value = 23124  # K132: 4668 -> 23124
$ws.Cells.Item(132, 13).Value = -20594  # M132: -2138 -> -20594
$ws.Cells.Item(134, 8).Value = 6259.968  # H134: 6435.4 -> 6259.968
$ws.Cells.Item(134, 9).Value = 2530.2104  # I134: 2615.389 -> 2530.2104
$ws.Cells.Item(134, 11).Value = 7590.6312  # K134: 7846.167 -> 7590.6312
$ws.Cells.Item(134, 13).Value = -5055.6312  # M134: -5311.167 -> -5055.6312

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(39, 8).Value = 1116.9412  # H39: 1078.3158 -> 1116.9412
$ws.Cells.Item(39, 9).Value = 829.3333  # I39: 747 -> 829.3333
$ws.Cells.Item(39, 10).Value = 1178.5714  # J39: 1166.6666 -> 1178.5714
$ws.Cells.Item(39, 11).Value = 2487.9999  # K39: 2241 -> 2487.9999
$ws.Cells.Item(39, 12).Value = 3535.7142  # L39: 3499.9998 -> 3535.7142
$ws.Cells.Item(39, 13).Value = -2193.9999  # M39: -1947 -> -2193.9999
$ws.Cells.Item(39, 14).Value = -4123.7142  # N39: -4087.9998 -> -4123.7142
$ws.Cells.Item(86, 8).Value = 866.6667  # H86: 849.75 -> 866.6667
$ws.Cells.Item(86, 9).Value = 866.6667  # I86: 849.75 -> 866.6667
$ws.Cells.Item(86, 11).Value = 2600.0001  # K86: 2549.25 -> 2600.0001
$ws.Cells.Item(86, 13).Value = -1414.0001  # M86: -1363.25 -> -1414.0001
$ws.Cells.Item(89, 8).Value = 866.6667  # H89: 849.75 -> 866.6667
$ws.Cells.Item(89, 9).Value = 866.6667  # I89: 849.75 -> 866.6667
$ws.Cells.Item(89, 11).Value = 7800.0003  # K89: 7647.75 -> 7800.0003
$ws.Cells.Item(89, 13).Value = -1872.0003  # M89: -1719.75 -> -1872.0003
$ws.Cells.Item(112, 8).Value = 7030.2  # H112: 6691.6665 -> 7030.2
$ws.Cells.Item(112, 9).Value = 4993  # I112: 4995 -> 4993
$ws.Cells.Item(112, 11).Value = 14979  # K112: 14985 -> 14979
$ws.Cells.Item(112, 13).Value = -13871  # M112: -13877 -> -13871
$ws.Cells.Item(128, 8).Value = 417989.8  # H128: 417995.8 -> 417989.8
$ws.Cells.Item(128, 9).Value = 417989.8  # I128: 417995.8 -> 417989.8
$ws.Cells.Item(128, 11).Value = 1253969.4  # K128: 1253987.4 -> 1253969.4
$ws.Cells.Item(128, 13).Value = -1248989.4  # M128: -1249007.4 -> -1248989.4
$ws.Cells.Item(138, 8).Value = 1999.5  # H138: 2000 -> 1999.5
$ws.Cells.Item(138, 9).Value = 1999.5  # I138: 2000 -> 1999.5
$ws.Cells.Item(138, 11).Value = 5998.5  # K138: 6000 -> 5998.5
$ws.Cells.Item(138, 13).Value = -858.5  # M138: -860 -> -858.5

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(102, 8).Value = 2619.7273  # H102: 2725.476 -> 2619.7273
$ws.Cells.Item(102, 9).Value = 2116.0667  # I102: 2238.7144 -> 2116.0667
$ws.Cells.Item(102, 11).Value = 2116.0667  # K102: 2238.7144 -> 2116.0667
$ws.Cells.Item(102, 13).Value = -494.0666999999999  # M102: -616.7143999999998 -> -494.0666999999999
$ws.Cells.Item(103, 8).Value = 0  # H103: 45999.5 -> 0
$ws.Cells.Item(103, 10).Value = 0  # J103: 45999.5 -> 0
$ws.Cells.Item(103, 12).Value = 0  # L103: 45999.5 -> 0
$ws.Cells.Item(103, 14).ClearContents()  # N103: remove (was -48343.5)
$ws.Cells.Item(135, 8).Value = 100780  # H135: 0 -> 100780
$ws.Cells.Item(135, 10).Value = 100780  # J135: 0 -> 100780
$ws.Cells.Item(135, 12).Value = 100780  # L135: 0 -> 100780
$ws.Cells.Item(135, 14).Value = -110920  # N135: None -> -110920

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(22, 8).Value = 2009.6364  # H22: 2278.25 -> 2009.6364
$ws.Cells.Item(22, 10).Value = 1856.1428  # J22: 0 -> 1856.1428
$ws.Cells.Item(22, 12).Value = 1856.1428  # L22: 0 -> 1856.1428
$ws.Cells.Item(22, 14).Value = -2446.1428  # N22: None -> -2446.1428
$ws.Cells.Item(27, 8).Value = 2009.6364  # H27: 2278.25 -> 2009.6364
$ws.Cells.Item(27, 10).Value = 1856.1428  # J27: 0 -> 1856.1428
$ws.Cells.Item(27, 12).Value = 1856.1428  # L27: 0 -> 1856.1428
$ws.Cells.Item(27, 14).Value = -2070.1428  # N27: None -> -2070.1428
$ws.Cells.Item(46, 8).Value = 3709.5  # H46: 3382.9167 -> 3709.5
$ws.Cells.Item(46, 9).Value = 2349  # I46: 2066 -> 2349
$ws.Cells.Item(46, 10).Value = 4049.625  # J46: 3821.889 -> 4049.625
$ws.Cells.Item(46, 11).Value = 2349  # K46: 2066 -> 2349
$ws.Cells.Item(46, 12).Value = 4049.625  # L46: 3821.889 -> 4049.625
$ws.Cells.Item(46, 13).Value = -2161  # M46: -1878 -> -2161
$ws.Cells.Item(46, 14).Value = -4425.625  # N46: -4197.889 -> -4425.625
$ws.Cells.Item(68, 8).Value = 12500  # H68: 15000 -> 12500
$ws.Cells.Item(68, 10).Value = 10000  # J68: 0 -> 10000
$ws.Cells.Item(68, 12).Value = 10000  # L68: 0 -> 10000
$ws.Cells.Item(68, 14).Value = -11498  # N68: None -> -11498
$ws.Cells.Item(71, 8).Value = 12500  # H71: 15000 -> 12500
$ws.Cells.Item(71, 10).Value = 10000  # J71: 0 -> 10000
$ws.Cells.Item(71, 12).Value = 50000  # L71: 0 -> 50000
$ws.Cells.Item(71, 14).Value = -57488  # N71: None -> -57488
$ws.Cells.Item(100, 8).Value = 3116.3333  # H100: 3193.5881 -> 3116.3333
$ws.Cells.Item(100, 9).Value = 2682.4167  # I100: 2762.3635 -> 2682.4167
$ws.Cells.Item(100, 11).Value = 2682.4167  # K100: 2762.3635 -> 2682.4167
$ws.Cells.Item(100, 13).Value = -2141.4167  # M100: -2221.3635 -> -2141.4167
$ws.Cells.Item(122, 8).Value = 4380.125  # H122: 4720.143 -> 4380.125
$ws.Cells.Item(122, 9).Value = 3591.75  # I122: 4122.3335 -> 3591.75
$ws.Cells.Item(122, 11).Value = 10775.25  # K122: 12367.0005 -> 10775.25
$ws.Cells.Item(122, 13).Value = -8325.25  # M122: -9917.000499999998 -> -8325.25

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(107, 8).Value = 2065.4167  # H107: 1982.6923 -> 2065.4167
$ws.Cells.Item(107, 10).Value = 4171  # J107: 3375.75 -> 4171
$ws.Cells.Item(107, 12).Value = 12513  # L107: 10127.25 -> 12513
$ws.Cells.Item(107, 14).Value = -16353  # N107: -13967.25 -> -16353
$ws.Cells.Item(122, 8).Value = 5293  # H122: 5296.6665 -> 5293
$ws.Cells.Item(122, 9).Value = 4507.9165  # I122: 4516.1665 -> 4507.9165
$ws.Cells.Item(122, 11).Value = 13523.7495  # K122: 13548.4995 -> 13523.7495
$ws.Cells.Item(122, 13).Value = -11073.7495  # M122: -11098.4995 -> -11073.7495
$ws.Cells.Item(136, 8).Value = 5564.1035  # H136: 5584.793 -> 5564.1035
$ws.Cells.Item(136, 9).Value = 4230.0527  # I136: 4261.6313 -> 4230.0527
$ws.Cells.Item(136, 11).Value = 12690.1581  # K136: 12784.8939 -> 12690.1581
$ws.Cells.Item(136, 13).Value = -10140.1581  # M136: -10234.8939 -> -10140.1581
